# "Hora de tareas modificada"
# - Header row: "Horas" -> " Horas estimadas", new "Horas reales" column (C)
#   inserted, names column moved from C to D with a new "Autor" header.
# - Several duration cells switch from text to real numbers.
# - Four new task rows (11-14) are appended.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row ---------------------------------------------------------
$ws.Range("B1").Value = " Horas estimadas"
$ws.Range("C1").Value = "Horas reales"
$ws.Range("D1").Value = "Autor"

# --- Existing rows: move the "Autor" names from C to D, fill in the new
#     "Horas reales" column, and normalise a few duration values to numbers
$ws.Range("D2").Value = "Clara Adolfo"
$ws.Range("C2").Value = 2

$ws.Range("D3").Value = "Clara Adolfo Rafael "
$ws.Range("C3").Value = 2

$ws.Range("D4").Value = "Todos"
$ws.Range("C4").Value = 0.5

$ws.Range("D5").Value = "Ángel"
$ws.Range("C5").ClearContents()

$ws.Range("D6").Value = "Ángel"
$ws.Range("C6").ClearContents()

$ws.Range("B7").Value = 0.5
$ws.Range("D7").Value = "Jose"
$ws.Range("C7").ClearContents()

$ws.Range("B8").Value = 0.5
$ws.Range("D8").Value = "Jose"
$ws.Range("C8").ClearContents()

$ws.Range("B9").Value = 0.0833
$ws.Range("D9").Value = "Todos"

$ws.Range("B10").Value = 0.5
$ws.Range("D10").Value = "Todos"

# --- New task rows -------------------------------------------------------
$ws.Range("A11").Value = "Buscar la sintaxis de java para los test"
$ws.Range("B11").Value = 0.5
$ws.Range("D11").Value = "Clara"

$ws.Range("A12").Value = "Pensar qué probar"
$ws.Range("B12").Value = 0.5
$ws.Range("D12").Value = "Clara"

$ws.Range("A13").Value = "Hacer test"
$ws.Range("B13").Value = 4
$ws.Range("D13").Value = "Clara"

$ws.Range("A14").Value = "Modificar la hoja de tareas"
$ws.Range("B14").Value = 0.5
$ws.Range("C14").Value = 0.5
$ws.Range("D14").Value = "Clara"

# --- Column widths -----------------------------------------------------
# (the runtime quantises ColumnWidth to the nearest 1/6 of a character, so
# these inputs are chosen to land on the closest achievable value to the
# authored widths of 40.85546875 / 17.28515625 / 19.5703125)
$ws.Columns.Item(1).ColumnWidth = 40.0029
$ws.Columns.Item(2).ColumnWidth = 16.5027
$ws.Columns.Item(4).ColumnWidth = 18.6678

# --- Selection -------------------------------------------------------------
$ws.Range("D1").Select()
